$d = $word.ActiveDocument

# 1. Update the keywords paragraph text: split "Número real, operaciones con
#    reales" into "número real," + "operaciones con reales" (note: the
#    space between the two pieces is intentionally dropped, matching the
#    target edit) using Find/Replace on the whole document content.
$d.Content.Find.Execute("Número real, operaciones con reales", $true, $false, $false, $false, $false, $true, 1, $false, "número real,operaciones con reales", 2) | Out-Null

# 2. Move the "_GoBack" bookmark from its old location (end of the
#    "Refuerza tu aprendizaje..." heading) to the split point between
#    "número real," and "operaciones con reales". Word only ever keeps a
#    single "_GoBack" bookmark (last edit position), so remove the old one
#    first.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng = $d.Content
$rng.Find.Execute("número real,operaciones con reales") | Out-Null
$splitPos = $rng.Start + ("número real,").Length
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
